{"js": "// The edit replaces every arithmetic-expression cell in the document's\n// (single) table with a new expression, preserving cell order, run\n// formatting (font/size) and paragraph properties. Word's Office.js API\n// exposes a table's cell text as a 2-D `values` array (row-major), so we\n// just overwrite that array in one shot instead of touching 100\n// individual cells/ranges.\nconst table = context.document.body.tables.getFirst();\n\nconst newValues = [\n  [\"44+38=\", \"84-36=\", \"70-13=\", \"70-22=\", \"75-16=\"],\n  [\"60-3=\", \"36+59=\", \"25+7=\", \"3+49=\", \"63-44=\"],\n  [\"17+55=\", \"41-19=\", \"25+49=\", \"71-49=\", \"9+54=\"],\n  [\"24+17=\", \"92-24=\", \"52-9=\", \"49+38=\", \"18+37=\"],\n  [\"71-64=\", \"91-87=\", \"5+66=\", \"68+9=\", \"53-28=\"],\n  [\"62-56=\", \"41-38=\", \"34+57=\", \"84-29=\", \"58+3=\"],\n  [\"50-14=\", \"46+6=\", \"65+29=\", \"58+27=\", \"34+9=\"],\n  [\"10-6=\", \"16+58=\", \"28+57=\", \"46+39=\", \"63-6=\"],\n  [\"72-28=\", \"28+33=\", \"35+8=\", \"7+48=\", \"49+27=\"],\n  [\"63-37=\", \"18+78=\", \"55+6=\", \"54-5=\", \"80-6=\"],\n  [\"48+34=\", \"84-39=\", \"67+7=\", \"94-68=\", \"60-2=\"],\n  [\"35+56=\", \"74-37=\", \"8+75=\", \"66-17=\", \"29+15=\"],\n  [\"16+46=\", \"21-19=\", \"28-19=\", \"55+8=\", \"39+4=\"],\n  [\"63-19=\", \"84-66=\", \"34+19=\", \"85-47=\", \"61-49=\"],\n  [\"92-76=\", \"92-55=\", \"33+59=\", \"36+49=\", \"82-38=\"],\n  [\"90-38=\", \"9+19=\", \"75-46=\", \"55+36=\", \"60-26=\"],\n  [\"15+6=\", \"40-7=\", \"75-9=\", \"27+8=\", \"18+9=\"],\n  [\"57+16=\", \"45+16=\", \"80-11=\", \"38-29=\", \"73-64=\"],\n  [\"49+18=\", \"44-6=\", \"66+8=\", \"90-4=\", \"28+15=\"],\n  [\"10-9=\", \"55+16=\", \"90-74=\", \"76+19=\", \"17+57=\"]\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# The edit replaces every arithmetic-expression cell in the document's\n# (single) table with a new expression, in cell order (row-major),\n# preserving each cell's run formatting (font/size) and paragraph\n# properties untouched. We walk Table 1's Rows/Cells in order and\n# assign the new text straight onto Cell.Range.Text (Word re-renders\n# the end-of-cell marker automatically).\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$newValues = @(\n    \"44+38=\", \"84-36=\", \"70-13=\", \"70-22=\", \"75-16=\",\n    \"60-3=\", \"36+59=\", \"25+7=\", \"3+49=\", \"63-44=\",\n    \"17+55=\", \"41-19=\", \"25+49=\", \"71-49=\", \"9+54=\",\n    \"24+17=\", \"92-24=\", \"52-9=\", \"49+38=\", \"18+37=\",\n    \"71-64=\", \"91-87=\", \"5+66=\", \"68+9=\", \"53-28=\",\n    \"62-56=\", \"41-38=\", \"34+57=\", \"84-29=\", \"58+3=\",\n    \"50-14=\", \"46+6=\", \"65+29=\", \"58+27=\", \"34+9=\",\n    \"10-6=\", \"16+58=\", \"28+57=\", \"46+39=\", \"63-6=\",\n    \"72-28=\", \"28+33=\", \"35+8=\", \"7+48=\", \"49+27=\",\n    \"63-37=\", \"18+78=\", \"55+6=\", \"54-5=\", \"80-6=\",\n    \"48+34=\", \"84-39=\", \"67+7=\", \"94-68=\", \"60-2=\",\n    \"35+56=\", \"74-37=\", \"8+75=\", \"66-17=\", \"29+15=\",\n    \"16+46=\", \"21-19=\", \"28-19=\", \"55+8=\", \"39+4=\",\n    \"63-19=\", \"84-66=\", \"34+19=\", \"85-47=\", \"61-49=\",\n    \"92-76=\", \"92-55=\", \"33+59=\", \"36+49=\", \"82-38=\",\n    \"90-38=\", \"9+19=\", \"75-46=\", \"55+36=\", \"60-26=\",\n    \"15+6=\", \"40-7=\", \"75-9=\", \"27+8=\", \"18+9=\",\n    \"57+16=\", \"45+16=\", \"80-11=\", \"38-29=\", \"73-64=\",\n    \"49+18=\", \"44-6=\", \"66+8=\", \"90-4=\", \"28+15=\",\n    \"10-9=\", \"55+16=\", \"90-74=\", \"76+19=\", \"17+57=\"\n)\n\n$i = 0\nforeach ($row in $table.Rows) {\n    foreach ($cell in $row.Cells) {\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n"}
